# Applies the "scheduled runner" market data refresh to the Aegis_Profits sheets.
# For each affected leve row, the currentAveragePrice* (H:L) and LeveProfit* (M:N)
# columns are updated to the latest fetched values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2785.0908
$ws.Range("I43").Value = 2787.24
$ws.Range("J43").Value = 2778.375
$ws.Range("K43").Value = 2787.24
$ws.Range("L43").Value = 2778.375
$ws.Range("M43").Value = -2718.24
$ws.Range("N43").Value = -2916.375
$ws.Range("H58").Value = 1227100.5
$ws.Range("I58").Value = 2451238.5
$ws.Range("J58").Value = 2962.5
$ws.Range("K58").Value = 7353715.5
$ws.Range("L58").Value = 8887.5
$ws.Range("M58").Value = -7353565.5
$ws.Range("N58").Value = -9187.5
$ws.Range("H86").Value = 4167.6665
$ws.Range("I86").Value = 3000.6667
$ws.Range("K86").Value = 3000.6667
$ws.Range("M86").Value = -1877.6667
$ws.Range("H89").Value = 4167.6665
$ws.Range("I89").Value = 3000.6667
$ws.Range("K89").Value = 15003.3335
$ws.Range("M89").Value = -9387.333500000001
$ws.Range("H132").Value = 3792365.2
$ws.Range("I132").Value = 4469095
$ws.Range("J132").Value = 2680.6
$ws.Range("K132").Value = 13407285
$ws.Range("L132").Value = 8041.799999999999
$ws.Range("M132").Value = -13404755
$ws.Range("N132").Value = -13101.8
$ws.Range("H135").Value = 1323.3334
$ws.Range("I135").Value = 557.3939
$ws.Range("J135").Value = 3008.4
$ws.Range("K135").Value = 5016.5451
$ws.Range("L135").Value = 27075.6
$ws.Range("M135").Value = -2481.5451
$ws.Range("N135").Value = -32145.6
$ws.Range("H137").Value = 1086.2
$ws.Range("I137").Value = 1035.8846
$ws.Range("J137").Value = 1413.25
$ws.Range("K137").Value = 3107.6538
$ws.Range("L137").Value = 4239.75
$ws.Range("M137").Value = -557.6538
$ws.Range("N137").Value = -9339.75
$ws.Range("H138").Value = 3458.8025
$ws.Range("I138").Value = 3087
$ws.Range("J138").Value = 3529.8823
$ws.Range("K138").Value = 9261
$ws.Range("L138").Value = 10589.6469
$ws.Range("M138").Value = -4121
$ws.Range("N138").Value = -20869.6469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29833.875
$ws.Range("I32").Value = 13206.167
$ws.Range("J32").Value = 65464.68
$ws.Range("K32").Value = 13206.167
$ws.Range("L32").Value = 65464.68
$ws.Range("M32").Value = -12919.167
$ws.Range("N32").Value = -66038.67999999999
$ws.Range("H44").Value = 13008.167
$ws.Range("J44").Value = 13009.8
$ws.Range("L44").Value = 13009.8
$ws.Range("N44").Value = -13985.8
$ws.Range("H55").Value = 12372
$ws.Range("J55").Value = 12372
$ws.Range("L55").Value = 12372
$ws.Range("N55").Value = -13002
$ws.Range("H61").Value = 2015.3334
$ws.Range("I61").Value = 1729.2
$ws.Range("K61").Value = 1729.2
$ws.Range("M61").Value = -1517.2
$ws.Range("H132").Value = 18323.457
$ws.Range("I132").Value = 20174.58
$ws.Range("K132").Value = 60523.74000000001
$ws.Range("M132").Value = -57993.74000000001
$ws.Range("H136").Value = 2015.3334
$ws.Range("I136").Value = 1729.2
$ws.Range("K136").Value = 5187.6
$ws.Range("M136").Value = -2637.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 87834.30499999999
$ws.Range("I86").Value = 103476.91
$ws.Range("K86").Value = 103476.91
$ws.Range("M86").Value = -102353.91
$ws.Range("H89").Value = 87834.30499999999
$ws.Range("I89").Value = 103476.91
$ws.Range("K89").Value = 517384.55
$ws.Range("M89").Value = -511768.55
$ws.Range("H117").Value = 36000
$ws.Range("J117").Value = 36000
$ws.Range("L117").Value = 36000
$ws.Range("N117").Value = -45178
$ws.Range("H134").Value = 3761.302
$ws.Range("I134").Value = 3882.8096
$ws.Range("J134").Value = 3297.3635
$ws.Range("K134").Value = 11648.4288
$ws.Range("L134").Value = 9892.0905
$ws.Range("M134").Value = -9113.4288
$ws.Range("N134").Value = -14962.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23468.188
$ws.Range("I31").Value = 1193.1025
$ws.Range("J31").Value = 52425.8
$ws.Range("K31").Value = 1193.1025
$ws.Range("L31").Value = 52425.8
$ws.Range("M31").Value = -898.1025
$ws.Range("N31").Value = -53015.8
$ws.Range("H34").Value = 23468.188
$ws.Range("I34").Value = 1193.1025
$ws.Range("J34").Value = 52425.8
$ws.Range("K34").Value = 1193.1025
$ws.Range("L34").Value = 52425.8
$ws.Range("M34").Value = -991.1025
$ws.Range("N34").Value = -52829.8
$ws.Range("H58").Value = 2163.4075
$ws.Range("I58").Value = 2096.5334
$ws.Range("J58").Value = 2247
$ws.Range("K58").Value = 2096.5334
$ws.Range("L58").Value = 2247
$ws.Range("M58").Value = -1893.5334
$ws.Range("N58").Value = -2653
$ws.Range("H105").Value = 1318.6666
$ws.Range("I105").Value = 1338.2222
$ws.Range("J105").Value = 1260
$ws.Range("K105").Value = 1338.2222
$ws.Range("L105").Value = 1260
$ws.Range("M105").Value = 408.7778000000001
$ws.Range("N105").Value = -4754
$ws.Range("H136").Value = 2163.4075
$ws.Range("I136").Value = 2096.5334
$ws.Range("J136").Value = 2247
$ws.Range("K136").Value = 6289.600199999999
$ws.Range("L136").Value = 6741
$ws.Range("M136").Value = -3739.600199999999
$ws.Range("N136").Value = -11841

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 5200
$ws.Range("J39").Value = 6766.6665
$ws.Range("L39").Value = 20299.9995
$ws.Range("N39").Value = -20887.9995
$ws.Range("H122").Value = 885.4167
$ws.Range("I122").Value = 470
$ws.Range("J122").Value = 1300.8334
$ws.Range("K122").Value = 4230
$ws.Range("L122").Value = 11707.5006
$ws.Range("M122").Value = -1780
$ws.Range("N122").Value = -16607.5006
$ws.Range("H133").Value = 3999.75
$ws.Range("I133").Value = 999
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 2997
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = 2063
$ws.Range("N133").Value = -25120
$ws.Range("H136").Value = 2916.25
$ws.Range("I136").Value = 2470
$ws.Range("K136").Value = 7410
$ws.Range("M136").Value = -2310
$ws.Range("H137").Value = 39886.2
$ws.Range("J137").Value = 8755.053
$ws.Range("L137").Value = 26265.159
$ws.Range("N137").Value = -36465.159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 100223620
$ws.Range("I80").Value = 167034800
$ws.Range("J80").Value = 6825
$ws.Range("K80").Value = 167034800
$ws.Range("L80").Value = 6825
$ws.Range("M80").Value = -167033802
$ws.Range("N80").Value = -8821
$ws.Range("H83").Value = 100223620
$ws.Range("I83").Value = 167034800
$ws.Range("J83").Value = 6825
$ws.Range("K83").Value = 835174000
$ws.Range("L83").Value = 34125
$ws.Range("M83").Value = -835169008
$ws.Range("N83").Value = -44109
$ws.Range("H97").Value = 47620588
$ws.Range("I97").Value = 62501756
$ws.Range("J97").Value = 844.4
$ws.Range("K97").Value = 62501756
$ws.Range("L97").Value = 844.4
$ws.Range("M97").Value = -62501260
$ws.Range("N97").Value = -1836.4
$ws.Range("H113").Value = 1208.2
$ws.Range("I113").Value = 812.3333
$ws.Range("K113").Value = 812.3333
$ws.Range("M113").Value = 1357.6667
$ws.Range("H122").Value = 5879
$ws.Range("I122").Value = 6473.75
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 19421.25
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -16971.25
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 2562.4583
$ws.Range("I132").Value = 1996.6207
$ws.Range("J132").Value = 3426.1052
$ws.Range("K132").Value = 5989.8621
$ws.Range("L132").Value = 10278.3156
$ws.Range("M132").Value = -3459.8621
$ws.Range("N132").Value = -15338.3156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3425.7144
$ws.Range("I7").Value = 1890
$ws.Range("J7").Value = 4040
$ws.Range("K7").Value = 1890
$ws.Range("L7").Value = 4040
$ws.Range("M7").Value = -1778
$ws.Range("N7").Value = -4264
$ws.Range("H40").Value = 504999.5
$ws.Range("I40").Value = 504999.5
$ws.Range("K40").Value = 504999.5
$ws.Range("M40").Value = -504863.5
$ws.Range("H82").Value = 2370.25
$ws.Range("I82").Value = 1993
$ws.Range("K82").Value = 1993
$ws.Range("M82").Value = -1632
$ws.Range("H85").Value = 2370.25
$ws.Range("I85").Value = 1993
$ws.Range("K85").Value = 1993
$ws.Range("M85").Value = -745
$ws.Range("H93").Value = 3435.7896
$ws.Range("I93").Value = 3492.0667
$ws.Range("J93").Value = 3224.75
$ws.Range("K93").Value = 3492.0667
$ws.Range("L93").Value = 3224.75
$ws.Range("M93").Value = -2244.0667
$ws.Range("N93").Value = -5720.75
$ws.Range("H122").Value = 4500.6665
$ws.Range("I122").Value = 4000.8
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 12002.4
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -9552.400000000001
$ws.Range("N122").Value = -25900
$ws.Range("H123").Value = 27350
$ws.Range("J123").Value = 27350
$ws.Range("L123").Value = 27350
$ws.Range("N123").Value = -37150
$ws.Range("H126").Value = 3425.7144
$ws.Range("I126").Value = 1890
$ws.Range("J126").Value = 4040
$ws.Range("K126").Value = 5670
$ws.Range("L126").Value = 12120
$ws.Range("M126").Value = -3200
$ws.Range("N126").Value = -17060
$ws.Range("H132").Value = 4473.1924
$ws.Range("I132").Value = 6389.1665
$ws.Range("J132").Value = 2830.9285
$ws.Range("K132").Value = 19167.4995
$ws.Range("L132").Value = 8492.7855
$ws.Range("M132").Value = -16637.4995
$ws.Range("N132").Value = -13552.7855
$ws.Range("H136").Value = 2786.238
$ws.Range("I136").Value = 2079.5789
$ws.Range("J136").Value = 9499.5
$ws.Range("K136").Value = 6238.736699999999
$ws.Range("L136").Value = 28498.5
$ws.Range("M136").Value = -3688.736699999999
$ws.Range("N136").Value = -33598.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 23133.125
$ws.Range("I132").Value = 2439.1943
$ws.Range("J132").Value = 85214.914
$ws.Range("K132").Value = 7317.5829
$ws.Range("L132").Value = 255644.742
$ws.Range("M132").Value = -4787.5829
$ws.Range("N132").Value = -260704.742
$ws.Range("H136").Value = 2885.1077
$ws.Range("I136").Value = 3424.4
$ws.Range("J136").Value = 2255.9333
$ws.Range("K136").Value = 10273.2
$ws.Range("L136").Value = 6767.7999
$ws.Range("M136").Value = -7723.200000000001
$ws.Range("N136").Value = -11867.7999
$ws.Range("H139").Value = 64625.8
$ws.Range("J139").Value = 64625.8
$ws.Range("L139").Value = 64625.8
$ws.Range("N139").Value = -74905.8
